# Apply the "UFC score" sheet addition to the workbook.

$wb = $excel.ActiveWorkbook

# --- 1) Restore the UFC sheet's selection (it will lose "active" status
#        once the new sheet becomes the active tab, which matches the
#        target: tabSelected moves off of UFC and the selection becomes B10).
$wsUfc = $wb.Worksheets.Item("UFC")
$wsUfc.Activate()
$wsUfc.Range("B10").Select()

# --- 2) Add the new "UFC score" sheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsScore = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsScore.Name = "UFC score"

# --- 3) Header row.
$headers = @("Catégorie", "Poids simple", "Poids moyen", "Poids difficile", "Score Simple", "Score Moyen", "Score difficile", "Total")
for ($c = 1; $c -le $headers.Length; $c++) {
    $wsScore.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# --- 4) Data rows (A:G). Column H (Total) is only populated on row 2.
$data = @(
    @("Entrées ext",  1, 3, 5, 2, 0, 0),
    @("Sorties ext",  2, 4, 6, 1, 0, 0),
    @("Requetes ext", 1, 3, 5, 0, 0, 0),
    @("Fichiers int", 3, 9, 4, 0, 1, 0),
    @("Fichiers int", 3, 9, 4, 0, 1, 0),
    @("Fichiers ext", 4, 6, 9, 0, 0, 0)
)

$row = 2
foreach ($entry in $data) {
    for ($c = 1; $c -le 7; $c++) {
        $wsScore.Cells.Item($row, $c).Value = $entry[$c - 1]
    }
    $row++
}

$wsScore.Cells.Item(2, 8).Value = 13

# --- 5) Final selection/activation on the new sheet (becomes the active tab).
$wsScore.Activate()
$wsScore.Range("D13").Select()
